# Auto-generated: apply scheduled Kraken_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 453.77777
$ws.Range("I19").Value = 448
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 448
$ws.Range("L19").Value = 500
$ws.Range("M19").Value = -273
$ws.Range("N19").Value = -850
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H70").Value = 1419.8
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 1524.75
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 4574.25
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -5114.25
$ws.Range("H73").Value = 1419.8
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 1524.75
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 4574.25
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -6446.25
$ws.Range("H92").Value = 2486.2222
$ws.Range("I92").Value = 2230.6667
$ws.Range("K92").Value = 2230.6667
$ws.Range("M92").Value = -982.6667000000002
$ws.Range("H107").Value = 1271.2858
$ws.Range("I107").Value = 804.4545000000001
$ws.Range("K107").Value = 804.4545000000001
$ws.Range("M107").Value = 1115.5455

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 8443.429
$ws.Range("I122").Value = 9360.182000000001
$ws.Range("J122").Value = 5082
$ws.Range("K122").Value = 28080.546
$ws.Range("L122").Value = 15246
$ws.Range("M122").Value = -25630.546
$ws.Range("N122").Value = -20146
$ws.Range("H132").Value = 1946
$ws.Range("I132").Value = 1770.4166
$ws.Range("K132").Value = 5311.2498
$ws.Range("M132").Value = -2781.2498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 10025000
$ws.Range("I7").Value = 10025000
$ws.Range("K7").Value = 10025000
$ws.Range("M7").Value = -10024887
$ws.Range("H80").Value = 608.8
$ws.Range("I80").Value = 515.4286
$ws.Range("J80").Value = 826.6667
$ws.Range("K80").Value = 515.4286
$ws.Range("L80").Value = 826.6667
$ws.Range("M80").Value = 482.5714
$ws.Range("N80").Value = -2822.6667
$ws.Range("H83").Value = 608.8
$ws.Range("I83").Value = 515.4286
$ws.Range("J83").Value = 826.6667
$ws.Range("K83").Value = 2577.143
$ws.Range("L83").Value = 4133.3335
$ws.Range("M83").Value = 2414.857
$ws.Range("N83").Value = -14117.3335
$ws.Range("H134").Value = 8548.9
$ws.Range("J134").Value = 11798.4
$ws.Range("L134").Value = 35395.2
$ws.Range("N134").Value = -40465.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 4346.3335
$ws.Range("I38").Value = 3019.5
$ws.Range("J38").Value = 7000
$ws.Range("K38").Value = 3019.5
$ws.Range("L38").Value = 7000
$ws.Range("M38").Value = -2642.5
$ws.Range("N38").Value = -7754
$ws.Range("H46").Value = 4346.3335
$ws.Range("I46").Value = 3019.5
$ws.Range("J46").Value = 7000
$ws.Range("K46").Value = 3019.5
$ws.Range("L46").Value = 7000
$ws.Range("M46").Value = -2808.5
$ws.Range("N46").Value = -7422
$ws.Range("H127").Value = 75000
$ws.Range("J127").Value = 75000
$ws.Range("L127").Value = 75000
$ws.Range("N127").Value = -84920
$ws.Range("H132").Value = 865.44446
$ws.Range("I132").Value = 672
$ws.Range("J132").Value = 1832.6666
$ws.Range("K132").Value = 2016
$ws.Range("L132").Value = 5497.9998
$ws.Range("M132").Value = 514
$ws.Range("N132").Value = -10557.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2273
$ws.Range("I129").Value = 2019.6666
$ws.Range("K129").Value = 6058.9998
$ws.Range("M129").Value = -1058.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2841.389
$ws.Range("I132").Value = 2202.5
$ws.Range("J132").Value = 5077.5
$ws.Range("K132").Value = 6607.5
$ws.Range("L132").Value = 15232.5
$ws.Range("M132").Value = -4077.5
$ws.Range("N132").Value = -20292.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4045.1667
$ws.Range("I7").Value = 4045.1667
$ws.Range("K7").Value = 4045.1667
$ws.Range("M7").Value = -3933.1667
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H38").Value = 25000
$ws.Range("I38").Value = 20000
$ws.Range("J38").Value = 30000
$ws.Range("K38").Value = 20000
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = -19590
$ws.Range("N38").Value = -30820
$ws.Range("H88").Value = 3000
$ws.Range("I88").Value = 3000
$ws.Range("K88").Value = 3000
$ws.Range("M88").Value = -2572
$ws.Range("H91").Value = 3000
$ws.Range("I91").Value = 3000
$ws.Range("K91").Value = 3000
$ws.Range("M91").Value = -1518
$ws.Range("H114").Value = 49999.5
$ws.Range("J114").Value = 49999.5
$ws.Range("L114").Value = 49999.5
$ws.Range("N114").Value = -58677.5
$ws.Range("H122").Value = 5599
$ws.Range("I122").Value = 4665
$ws.Range("K122").Value = 13995
$ws.Range("M122").Value = -11545
$ws.Range("H126").Value = 4045.1667
$ws.Range("I126").Value = 4045.1667
$ws.Range("K126").Value = 12135.5001
$ws.Range("M126").Value = -9665.500100000001
$ws.Range("H132").Value = 3857.5
$ws.Range("I132").Value = 4217.6
$ws.Range("J132").Value = 3257.3333
$ws.Range("K132").Value = 12652.8
$ws.Range("L132").Value = 9771.999899999999
$ws.Range("M132").Value = -10122.8
$ws.Range("N132").Value = -14831.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H81").Value = 643.75
$ws.Range("I81").Value = 525
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 1050
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = 11
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 643.75
$ws.Range("I84").Value = 525
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 5250
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = 54
$ws.Range("N84").Value = -20608
$ws.Range("H95").Value = 27344
$ws.Range("J95").Value = 27344
$ws.Range("L95").Value = 27344
$ws.Range("N95").Value = -32836
$ws.Range("H107").Value = 4747.5
$ws.Range("I107").Value = 2495
$ws.Range("J107").Value = 7000
$ws.Range("K107").Value = 7485
$ws.Range("L107").Value = 21000
$ws.Range("M107").Value = -5565
$ws.Range("N107").Value = -24840
$ws.Range("H111").Value = 35000
$ws.Range("J111").Value = 35000
$ws.Range("L111").Value = 35000
$ws.Range("N111").Value = -43180
$ws.Range("H113").Value = 473.66666
$ws.Range("I113").Value = 458
$ws.Range("K113").Value = 1374
$ws.Range("M113").Value = 796
$ws.Range("H132").Value = 3050.6316
$ws.Range("I132").Value = 1433.2142
$ws.Range("K132").Value = 4299.642599999999
$ws.Range("M132").Value = -1769.642599999999
